$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Ref, $Text)
    $r = $Sheet.Range($Ref)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "308.47"
Set-TextValue $ws "E2" "1.08%"
Set-TextValue $ws "D3" "36.21"
Set-TextValue $ws "E3" "0.83%"
Set-TextValue $ws "D4" "5.045"
Set-TextValue $ws "E4" "0.73%"
Set-TextValue $ws "D5" "0.08153"
Set-TextValue $ws "E5" "1.12%"
Set-TextValue $ws "D6" "1.982"
Set-TextValue $ws "E6" "4.95%"
Set-TextValue $ws "B7" "KuCoinToken"
Set-TextValue $ws "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D7" "7.863"
Set-TextValue $ws "E7" "0.33%"
Set-TextValue $ws "B8" "MXToken"
Set-TextValue $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D8" "0.9273"
Set-TextValue $ws "E8" "-0.48%"
Set-TextValue $ws "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D9" "0.1473"
Set-TextValue $ws "E9" "14.60%"
Set-TextValue $ws "B10" "WazirX"
Set-TextValue $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1946"
Set-TextValue $ws "E10" "2.34%"
Set-TextValue $ws "B11" "MandalaExchangeToken"
Set-TextValue $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.09099"
Set-TextValue $ws "E11" "-1.48%"
Set-TextValue $ws "B12" "BitrueCoin"
Set-TextValue $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D12" "0.03524"
Set-TextValue $ws "E12" "0.37%"
Set-TextValue $ws "B13" "BitMartToken"
Set-TextValue $ws "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.09877"
Set-TextValue $ws "E13" "-0.19%"
Set-TextValue $ws "B14" "BitForexToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001407"
Set-TextValue $ws "E14" "-1.90%"
Set-TextValue $ws "B15" "TigerCash"
Set-TextValue $ws "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D15" "0.006500"
Set-TextValue $ws "E15" "1.86%"
Set-TextValue $ws "B16" "LEO"
Set-TextValue $ws "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D16" "3.849"
Set-TextValue $ws "E16" "5.19%"
Set-TextValue $ws "B17" "GateToken"
Set-TextValue $ws "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D17" "4.155"
Set-TextValue $ws "E17" "0.14%"
Set-TextValue $ws "D18" "3.411"
Set-TextValue $ws "E18" "7.82%"
Set-TextValue $ws "D19" "0.3450"
Set-TextValue $ws "E19" "-0.04%"
Set-TextValue $ws "D20" "0.1294"
Set-TextValue $ws "E20" "-3.78%"
Set-TextValue $ws "D21" "4.832"
Set-TextValue $ws "E21" "-7.28%"
Set-TextValue $ws "D22" "0.2345"
Set-TextValue $ws "E22" "-7.55%"
Set-TextValue $ws "D23" "0.04365"
Set-TextValue $ws "E23" "-1.08%"
Set-TextValue $ws "E24" "-0.01%"
Set-TextValue $ws "D25" "0.004160"
Set-TextValue $ws "E25" "-11.60%"
Set-TextValue $ws "E27" "0.05%"
Set-TextValue $ws "D39" "0.02143"
Set-TextValue $ws "E39" "9.83%"
Set-TextValue $ws "D40" "0.05112"
Set-TextValue $ws "E40" "-0.86%"
Set-TextValue $ws "D41" "0.007450"
Set-TextValue $ws "E41" "-1.27%"
Set-TextValue $ws "D42" "0.009994"
Set-TextValue $ws "E42" "-1.60%"
Set-TextValue $ws "E43" "-0.14%"
Set-TextValue $ws "E44" "-1.80%"
Set-TextValue $ws "E45" "-9.90%"
Set-TextValue $ws "D46" "0.00006283"
Set-TextValue $ws "E46" "-1.05%"
Set-TextValue $ws "E47" "-0.02%"
Set-TextValue $ws "E48" "-0.64%"
Set-TextValue $ws "E49" "-3.56%"
Set-TextValue $ws "D50" "0.00002103"
Set-TextValue $ws "E50" "-0.02%"
Set-TextValue $ws "D51" "0.0002003"
Set-TextValue $ws "E51" "-0.02%"
